# Added the main dB and can now store new entries into main dB
#
# Appends two new rows to the documentation table on Sheet1, describing
# the new "main database" feature and the "save form to dBase" feature.
# Excel takes care of extending the used range / dimension, growing the
# shared-string table, and carrying the column C style (wrap-text, s="1")
# onto the new cells automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = 12
$ws.Range("C14").Value = "Main database of all entries"

$ws.Range("B15").Value = 13
$ws.Range("C15").Value = "Save the form to the dBase. checks locally if it's a edit or a new entry"

# Move the active selection to the next empty row in column C, matching
# where the user would land after typing the two new rows.
$ws.Range("C16").Select() | Out-Null
